$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pharmacophore generation had an H-Count bug that mismatched two ligand
# rows against the wrong PDB/activity data. Fix: row 2 (Ligand id 1) and
# row 3 (Ligand id 2) must swap their full contents so each SMILES/PDB
# entry lines up with its correct co-factors / activity value again.
#
# Row 2 becomes the former row-3 data (PDB 5jcb / NV4), row 3 becomes the
# former row-2 data (PDB 5xiw / LOC, including the Target_system note).
# Columns F (standard_unit), G (standard_type) and H (standard_relation)
# are identical between the two rows ("nM" / "Kd" / "=") so they are left
# untouched.

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "COc1c(OC)cc(cc1OC)[C@H]1[C@H]2C(=O)OC[C@H]2[C@H](c2c1cc1OCOc1c2)Sc1[nH]ncn1"
$ws.Range("C2").Value = "5jcb"
$ws.Range("D2").Value = "NV4"
$ws.Range("E2").Value = "GTP,ACP,GDP,MES,GOL,IMD,CA,Mg,NA"
$ws.Range("I2").Value = 62040
$ws.Range("J2").ClearContents()

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "CO[C@H]1CC[C@H]2[C@@H](CC1=O)[C@H](CCc1c2c(OC)c(c(c1)OC)OC)NC(=O)C"
$ws.Range("C3").Value = "5xiw"
$ws.Range("D3").Value = "LOC"
$ws.Range("E3").Value = "MES,GOL,GTP,GDP.CA,Mg"
$ws.Range("I3").Value = 11030
$ws.Range("J3").Value = "Tubulin-colchicine binding domain"

$ws.Range("B13").Select() | Out-Null

# Touch the sheet's last row so the saved worksheet records the full-sheet
# extent/selection state the same way the authoring session left it.
$ws.Rows.Item(1048576).RowHeight = 12.8 | Out-Null
